$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 496 (existing rows 496-532 shift down to 498-534)
$ws.Rows("496:497").Insert()

# New row 496 data
$ws.Cells.Item(496, 1).Value = 7
$ws.Cells.Item(496, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(496, 3).Value = "Ñuble"
$ws.Cells.Item(496, 4).Value = 45013
$ws.Cells.Item(496, 5).Value = 16
$ws.Cells.Item(496, 6).Value = 100114001
$ws.Cells.Item(496, 7).Value = "Papa"
$ws.Cells.Item(496, 8).Value = "Patagonia"
$ws.Cells.Item(496, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(496, 10).Value = 60
$ws.Cells.Item(496, 11).Value = 12000
$ws.Cells.Item(496, 12).Value = 12000
$ws.Cells.Item(496, 13).Value = 12000
$ws.Cells.Item(496, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(496, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(496, 16).Value = 480
$ws.Cells.Item(496, 17).Value = 25
$ws.Cells.Item(496, 18).Value = "Hortaliza"

# New row 497 data
$ws.Cells.Item(497, 1).Value = 7
$ws.Cells.Item(497, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(497, 3).Value = "Ñuble"
$ws.Cells.Item(497, 4).Value = 45013
$ws.Cells.Item(497, 5).Value = 16
$ws.Cells.Item(497, 6).Value = 100114001
$ws.Cells.Item(497, 7).Value = "Papa"
$ws.Cells.Item(497, 8).Value = "Patagonia"
$ws.Cells.Item(497, 9).Value = "1a (cosecha)"
$ws.Cells.Item(497, 10).Value = 70
$ws.Cells.Item(497, 11).Value = 10000
$ws.Cells.Item(497, 12).Value = 10000
$ws.Cells.Item(497, 13).Value = 10000
$ws.Cells.Item(497, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(497, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(497, 16).Value = 400
$ws.Cells.Item(497, 17).Value = 25
$ws.Cells.Item(497, 18).Value = "Hortaliza"
